# Apply cryptos list update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.394.80"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "3.570.84"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'589.61"
$ws.Range("E5").Value = "  +2.34%  "
$ws.Range("D6").Value = "'187.41"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "3.560.02"
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").Value = "'0.199"
$ws.Range("E10").Value = "  +7.99%  "
$ws.Range("D11").Value = "'0.649"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "'54.87"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("E13").Value = "  +1.28%  "
$ws.Range("D14").Value = "'9.60"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "4.144.95"
$ws.Range("D16").Value = "'19.50"
$ws.Range("E16").Value = "  -0.85%  "
$ws.Range("D17").Value = "70.431.00"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").Value = "3.583.52"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").Value = "'12.48"
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("E20").Value = "  -0.62%  "
$ws.Range("D21").Value = "'559.04"
$ws.Range("E21").Value = "  +14.56%  "
$ws.Range("E22").Value = "  -1.09%  "
$ws.Range("D23").Value = "'18.00"
$ws.Range("E23").Value = "  -8.78%  "
$ws.Range("D24").Value = "'4.67"
$ws.Range("E24").Value = "  +6.91%  "
$ws.Range("E25").Value = "  +0.91%  "
$ws.Range("D26").Value = "'95.92"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").Value = "'11.46"
$ws.Range("E27").Value = "  +1.82%  "
$ws.Range("D28").Value = "'2.99"
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("D29").Value = "'9.17"
$ws.Range("E29").Value = "  -1.75%  "
$ws.Range("D30").Value = "'32.26"
$ws.Range("E30").Value = "  +1.44%  "
$ws.Range("D31").Value = "'7.35"
$ws.Range("E31").Value = "  -2.42%  "
$ws.Range("D32").Value = "'12.56"
$ws.Range("E32").Value = "  +4.14%  "
$ws.Range("D33").Value = "'65.12"
$ws.Range("E33").Value = "  -2.83%  "
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("D35").Value = "'554.44"
$ws.Range("E35").Value = "  -1.83%  "
$ws.Range("E36").Value = "  +2.20%  "
$ws.Range("D37").Value = "'0.416"
$ws.Range("E37").Value = "  +5.38%  "
$ws.Range("D38").Value = "'38.20"
$ws.Range("E38").Value = "  -1.12%  "
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").Value = "0.0₃0771"
$ws.Range("E40").Value = "  -3.70%  "
$ws.Range("E41").Value = "  -1.39%  "
$ws.Range("D42").Value = "3.367.48"
$ws.Range("E42").Value = "  +2.78%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'3.39"
$ws.Range("E43").Value = "  -4.67%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'3.08"
$ws.Range("E44").Value = "  -5.55%  "
$ws.Range("E45").Value = "  +3.88%  "
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("D47").Value = "'0.0448"
$ws.Range("E47").Value = "  +2.01%  "
$ws.Range("D48").Value = "'9.23"
$ws.Range("E48").Value = "  -3.93%  "
$ws.Range("D49").Value = "'0.136"
$ws.Range("E49").Value = "  +0.15%  "
$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("B51").Value = "OceanProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean"
$ws.Range("D51").Value = "'1.46"
$ws.Range("E51").Value = "  +17.27%  "
